$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (s="1") from an existing header cell onto the two
# new header cells so they reuse the same cellXf instead of creating a new one.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in data rows 2-13: column I is always 1, column J mirrors column H
for ($r = 2; $r -le 13; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}
